$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C9:D16").Clear()
$ws.Range("A25:B31").EntireRow.Delete()
$ws.Range("A25:B31").Select()
$excel.ActiveWindow.ScrollRow = 6
